$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 31.27132566666667
$ws.Range("H2").Value = 93.81397700000001
$ws.Range("I2").Value = 0.9493361071405608
$ws.Range("J2").Value = 0.9493361071405608
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 3660.138647085684
$ws.Range("R2").Value = 32941.24782377115
$ws.Range("S2").Value = 0.3080942788078155
$ws.Range("T2").Value = 0.3080942788078155
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 31.27132566666667
$ws.Range("H3").Value = 93.81397700000001
$ws.Range("I3").Value = 0.9493361071405608
$ws.Range("J3").Value = 0.9493361071405608
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 3176.542428682825
$ws.Range("R3").Value = 28588.88185814543
$ws.Range("S3").Value = 0.2673872885790031
$ws.Range("T3").Value = 0.2673872885790031
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 31.27132566666667
$ws.Range("H4").Value = 93.81397700000001
$ws.Range("I4").Value = 0.9493361071405608
$ws.Range("J4").Value = 0.9493361071405608
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 4441.365982633727
$ws.Range("R4").Value = 39972.29384370354
$ws.Range("S4").Value = 0.3738545397537422
$ws.Range("T4").Value = 0.3738545397537422
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.96805
$ws.Range("H5").Value = 2.90415
$ws.Range("I5").Value = 0.02938809912676722
$ws.Range("J5").Value = 0.02938809912676721
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 113.30498921215
$ws.Range("R5").Value = 1019.74490290935
$ws.Range("S5").Value = 0.009537512729043749
$ws.Range("T5").Value = 0.009537512729043747
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.96805
$ws.Range("H6").Value = 2.90415
$ws.Range("I6").Value = 0.02938809912676722
$ws.Range("J6").Value = 0.02938809912676721
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 98.33455514053334
$ws.Range("R6").Value = 885.0109962648
$ws.Range("S6").Value = 0.008277367818301871
$ws.Range("T6").Value = 0.008277367818301868
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.96805
$ws.Range("H7").Value = 2.90415
$ws.Range("I7").Value = 0.02938809912676722
$ws.Range("J7").Value = 0.02938809912676721
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 137.4890334141333
$ws.Range("R7").Value = 1237.4013007272
$ws.Range("S7").Value = 0.0115732185794216
$ws.Range("T7").Value = 0.0115732185794216
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.700829
$ws.Range("H8").Value = 2.102487
$ws.Range("I8").Value = 0.02127579373267201
$ws.Range("J8").Value = 0.02127579373267201
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 82.028224042727
$ws.Range("R8").Value = 738.254016384543
$ws.Range("S8").Value = 0.006904773005922215
$ws.Range("T8").Value = 0.006904773005922215
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.700829
$ws.Range("H9").Value = 2.102487
$ws.Range("I9").Value = 0.02127579373267201
$ws.Range("J9").Value = 0.02127579373267201
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 71.19023598428268
$ws.Range("R9").Value = 640.712123858544
$ws.Range("S9").Value = 0.005992479118571026
$ws.Range("T9").Value = 0.005992479118571025
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.700829
$ws.Range("H10").Value = 2.102487
$ws.Range("I10").Value = 0.02127579373267201
$ws.Range("J10").Value = 0.02127579373267201
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 99.53649274169068
$ws.Range("R10").Value = 895.828434675216
$ws.Range("S10").Value = 0.008378541608178769
$ws.Range("T10").Value = 0.008378541608178769
